$wb = $excel.ActiveWorkbook

# ALC row 17: One for the Road
$ws = $wb.Sheets.Item("ALC")
$ws.Range("H17").Value = 2922.0715
$ws.Range("I17").Value = 2031.6666
$ws.Range("J17").Value = 3164.9092
$ws.Range("K17").Value = 6094.9998
$ws.Range("L17").Value = 9494.7276
$ws.Range("M17").Value = -5926.9998
$ws.Range("N17").Value = -9830.7276

# ALC row 32: Automata for the People
$ws = $wb.Sheets.Item("ALC")
$ws.Range("H32").Value = 4049.4
$ws.Range("I32").Value = 3598
$ws.Range("J32").Value = 4275.1
$ws.Range("K32").Value = 3598
$ws.Range("L32").Value = 4275.1
$ws.Range("M32").Value = -3272
$ws.Range("N32").Value = -4927.1

# ALC row 34: Sophomore Slump
$ws = $wb.Sheets.Item("ALC")
$ws.Range("H34").Value = 10261.571
$ws.Range("I34").Value = 10261.571
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 10261.571
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -10058.571
$ws.Range("N34").Value = ""

# ALC row 36: You Put Your Left Hand In
$ws = $wb.Sheets.Item("ALC")
$ws.Range("H36").Value = 10261.571
$ws.Range("I36").Value = 10261.571
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 10261.571
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = -9546.571
$ws.Range("N36").Value = ""

# ALC row 64: Forged from the Void
$ws = $wb.Sheets.Item("ALC")
$ws.Range("H64").Value = 5247.5
$ws.Range("J64").Value = 5500
$ws.Range("L64").Value = 5500
$ws.Range("N64").Value = -5996

# ALC row 67: Dodging the Draft (L)
$ws = $wb.Sheets.Item("ALC")
$ws.Range("H67").Value = 5247.5
$ws.Range("J67").Value = 5500
$ws.Range("L67").Value = 5500
$ws.Range("N67").Value = -7216

# ALC row 137: Cutting Edge of Culinary Quality
$ws = $wb.Sheets.Item("ALC")
$ws.Range("H137").Value = 1389.5
$ws.Range("I137").Value = 1389.4667
$ws.Range("J137").Value = 1389.6
$ws.Range("K137").Value = 4168.4001
$ws.Range("L137").Value = 4168.799999999999
$ws.Range("M137").Value = -1618.4001
$ws.Range("N137").Value = -9268.799999999999

# ARM row 61: Dealing with the Tough Stuff
$ws = $wb.Sheets.Item("ARM")
$ws.Range("H61").Value = 4343.6665
$ws.Range("I61").Value = 1274.75
$ws.Range("K61").Value = 1274.75
$ws.Range("M61").Value = -1062.75

# ARM row 97: Ore for Me
$ws = $wb.Sheets.Item("ARM")
$ws.Range("H97").Value = 789.2857
$ws.Range("I97").Value = 789.2857
$ws.Range("K97").Value = 789.2857
$ws.Range("M97").Value = -293.2857

# ARM row 132: Don't Bore Me, Ore Me
$ws = $wb.Sheets.Item("ARM")
$ws.Range("H132").Value = 1710.375
$ws.Range("I132").Value = 1710.375
$ws.Range("K132").Value = 5131.125
$ws.Range("M132").Value = -2601.125

# ARM row 136: Metal with Mettle
$ws = $wb.Sheets.Item("ARM")
$ws.Range("H136").Value = 4343.6665
$ws.Range("I136").Value = 1274.75
$ws.Range("K136").Value = 3824.25
$ws.Range("M136").Value = -1274.25

# BSM row 5: Axe Me Anything
$ws = $wb.Sheets.Item("BSM")
$ws.Range("H5").Value = 918.25
$ws.Range("I5").Value = 924.6667
$ws.Range("K5").Value = 924.6667
$ws.Range("M5").Value = -811.6667

# BSM row 94: High Steal
$ws = $wb.Sheets.Item("BSM")
$ws.Range("H94").Value = 516
$ws.Range("I94").Value = 519.4
$ws.Range("J94").Value = 499
$ws.Range("K94").Value = 519.4
$ws.Range("L94").Value = 499
$ws.Range("M94").Value = -68.39999999999998
$ws.Range("N94").Value = -1401

# CUL row 7: It's Always Sunny in Vylbrand
$ws = $wb.Sheets.Item("CUL")
$ws.Range("H7").Value = 61
$ws.Range("I7").Value = 58.75
$ws.Range("J7").Value = 70
$ws.Range("K7").Value = 176.25
$ws.Range("L7").Value = 210
$ws.Range("M7").Value = -64.25
$ws.Range("N7").Value = -434

# CUL row 34: Fever Pitch
$ws = $wb.Sheets.Item("CUL")
$ws.Range("H34").Value = 664.3333
$ws.Range("I34").Value = 297
$ws.Range("J34").Value = 848
$ws.Range("K34").Value = 891
$ws.Range("L34").Value = 2544
$ws.Range("M34").Value = -807
$ws.Range("N34").Value = -2712

# CUL row 39: Bloody Good Tart, This
$ws = $wb.Sheets.Item("CUL")
$ws.Range("H39").Value = 4199.5
$ws.Range("J39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("N39").Value = ""

# CUL row 74: The Nutcracker's Sweets
$ws = $wb.Sheets.Item("CUL")
$ws.Range("H74").Value = 16197.6
$ws.Range("I74").Value = 9990
$ws.Range("J74").Value = 17749.5
$ws.Range("K74").Value = 29970
$ws.Range("L74").Value = 53248.5
$ws.Range("M74").Value = -28909
$ws.Range("N74").Value = -55370.5

# CUL row 75: Breakfast of Champions
$ws = $wb.Sheets.Item("CUL")
$ws.Range("H75").Value = 4498
$ws.Range("I75").Value = 1492
$ws.Range("J75").Value = 6001
$ws.Range("K75").Value = 4476
$ws.Range("L75").Value = 18003
$ws.Range("M75").Value = -3478
$ws.Range("N75").Value = -19999

# CUL row 77: Time for a Midnight Snack (L)
$ws = $wb.Sheets.Item("CUL")
$ws.Range("H77").Value = 16197.6
$ws.Range("I77").Value = 9990
$ws.Range("J77").Value = 17749.5
$ws.Range("K77").Value = 89910
$ws.Range("L77").Value = 159745.5
$ws.Range("M77").Value = -84606
$ws.Range("N77").Value = -170353.5

# CUL row 78: Emerald Soup for the Soul (L)
$ws = $wb.Sheets.Item("CUL")
$ws.Range("H78").Value = 4498
$ws.Range("I78").Value = 1492
$ws.Range("J78").Value = 6001
$ws.Range("K78").Value = 13428
$ws.Range("L78").Value = 54009
$ws.Range("M78").Value = -8436
$ws.Range("N78").Value = -63993

# CUL row 107: Slippery Service
$ws = $wb.Sheets.Item("CUL")
$ws.Range("H107").Value = 993.5
$ws.Range("J107").Value = 993.5
$ws.Range("L107").Value = 2980.5
$ws.Range("N107").Value = -6820.5

# GSM row 2: Copper and Robbers
$ws = $wb.Sheets.Item("GSM")
$ws.Range("H2").Value = 211.66667
$ws.Range("I2").Value = 308
$ws.Range("K2").Value = 308
$ws.Range("M2").Value = -195

# GSM row 97: If I'd a Koppranickel for Every Time...
$ws = $wb.Sheets.Item("GSM")
$ws.Range("H97").Value = 304.45456
$ws.Range("J97").Value = 389
$ws.Range("L97").Value = 389
$ws.Range("N97").Value = -1381

# GSM row 126: Gold Rush Order
$ws = $wb.Sheets.Item("GSM")
$ws.Range("H126").Value = 3189.6365
$ws.Range("I126").Value = 2935.75
$ws.Range("J126").Value = 3866.6667
$ws.Range("K126").Value = 8807.25
$ws.Range("L126").Value = 11600.0001
$ws.Range("M126").Value = -6337.25
$ws.Range("N126").Value = -16540.0001

# GSM row 132: On Board for Lar
$ws = $wb.Sheets.Item("GSM")
$ws.Range("H132").Value = 949.25
$ws.Range("I132").Value = 949.25
$ws.Range("K132").Value = 2847.75
$ws.Range("M132").Value = -317.75

# LTW row 7: Tan Before the Ban
$ws = $wb.Sheets.Item("LTW")
$ws.Range("H7").Value = 6968.7144
$ws.Range("I7").Value = 4193.4287
$ws.Range("J7").Value = 8356.357
$ws.Range("K7").Value = 4193.4287
$ws.Range("L7").Value = 8356.357
$ws.Range("M7").Value = -4081.4287
$ws.Range("N7").Value = -8580.357

# LTW row 22: Skin off Their Backs
$ws = $wb.Sheets.Item("LTW")
$ws.Range("H22").Value = 2784.8572
$ws.Range("I22").Value = 2658.8
$ws.Range("J22").Value = 3100
$ws.Range("K22").Value = 2658.8
$ws.Range("L22").Value = 3100
$ws.Range("M22").Value = -2363.8
$ws.Range("N22").Value = -3690

# LTW row 27: Fire and Hide
$ws = $wb.Sheets.Item("LTW")
$ws.Range("H27").Value = 2784.8572
$ws.Range("I27").Value = 2658.8
$ws.Range("J27").Value = 3100
$ws.Range("K27").Value = 2658.8
$ws.Range("L27").Value = 3100
$ws.Range("M27").Value = -2551.8
$ws.Range("N27").Value = -3314

# LTW row 122: Hell on Leather
$ws = $wb.Sheets.Item("LTW")
$ws.Range("H122").Value = 5176.8057
$ws.Range("I122").Value = 4140.8335
$ws.Range("K122").Value = 12422.5005
$ws.Range("M122").Value = -9972.500499999998

# LTW row 126: Battered Books
$ws = $wb.Sheets.Item("LTW")
$ws.Range("H126").Value = 6968.7144
$ws.Range("I126").Value = 4193.4287
$ws.Range("J126").Value = 8356.357
$ws.Range("K126").Value = 12580.2861
$ws.Range("L126").Value = 25069.071
$ws.Range("M126").Value = -10110.2861
$ws.Range("N126").Value = -30009.071

# WVR row 49: A Leg Up on the Cold
$ws = $wb.Sheets.Item("WVR")
$ws.Range("H49").Value = 42999
$ws.Range("J49").Value = 42999
$ws.Range("L49").Value = 42999
$ws.Range("N49").Value = -43459

# WVR row 81: Where the Dragonflies, the Net Catches
$ws = $wb.Sheets.Item("WVR")
$ws.Range("H81").Value = 1253575
$ws.Range("I81").Value = 5019.6
$ws.Range("J81").Value = 3334500.8
$ws.Range("K81").Value = 10039.2
$ws.Range("L81").Value = 6669001.6
$ws.Range("M81").Value = -8978.200000000001
$ws.Range("N81").Value = -6671123.6

# WVR row 84: To Kill a Dragon on Nameday (L)
$ws = $wb.Sheets.Item("WVR")
$ws.Range("H84").Value = 1253575
$ws.Range("I84").Value = 5019.6
$ws.Range("J84").Value = 3334500.8
$ws.Range("K84").Value = 50196
$ws.Range("L84").Value = 33345008
$ws.Range("M84").Value = -44892
$ws.Range("N84").Value = -33355616

# WVR row 122: Heavy Armoire
$ws = $wb.Sheets.Item("WVR")
$ws.Range("H122").Value = 1899.5
$ws.Range("I122").Value = 1899.5
$ws.Range("K122").Value = 5698.5
$ws.Range("M122").Value = -3248.5
